$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve the existing "(m3/s)" shared string by copying it to F1 before the
# old units row (row 2) is removed, so the shared-string table keeps it in
# its original relative slot once everything else is compacted away.
$ws.Range("F1").Value = $ws.Range("F2").Text

# Remove the old "units" row (Hiver/Eté/Année -> (m3/s)/(MW)/(GWh)); this
# shifts every data row up by one and drops the now-unused shared strings.
$ws.Rows.Item(2).Delete()

# Build the new single header row.
$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

# E1 used to carry the old header style; the new header text in A1:E1 is
# unstyled, so clear the leftover formatting.
$ws.Range("E1").ClearFormats()

# F1:K1 need a style that uses the 9pt Arial header font but the default
# (General) number format. Creating a scratch named style, applying it, and
# then deleting the named style leaves behind exactly that cell format
# without leaving extra cellStyle/cellStyleXfs entries around.
$tmpStyle = $wb.Styles.Add("TempHeaderStyle")
$tmpStyle.Font.Name = "Arial"
$tmpStyle.Font.Size = 9
$ws.Range("F1:K1").Style = "TempHeaderStyle"
$wb.Styles.Item("TempHeaderStyle").Delete()

# Match the new active selection recorded in the sheet view.
$null = $ws.Range("A2:K2").Select()
